# Applies the cryptocurrency price/volume updates described in the commit diff.
# D-column price values that look numeric are prefixed with a leading apostrophe
# so Excel stores them as text (matching the original inlineStr/text cells)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.443.80"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.868.06"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'236.27"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4830"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "'0.2800"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "'0.06507"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "1.837.20"
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("D11").Value = "'0.07448"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'16.25"
$ws.Range("D13").Value = "'5.078"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "'87.14"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "'0.6424"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "30.418.44"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'13.00"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "'0.000007484"
$ws.Range("D20").Value = "'229.89"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").Value = "2.098.33"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'5.148"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'6.094"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "'169.35"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "'9.330"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").Value = "'18.34"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "'1.910"
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("D29").Value = "'0.1041"
$ws.Range("E29").Value = "  +11.27%  "
$ws.Range("E30").Value = "  -5.05%  "
$ws.Range("D31").Value = "'4.281"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "'3.991"
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").Value = "'0.04982"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "'1.179"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "'0.7428"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Value = "'0.9996"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'2.712"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'0.01933"
$ws.Range("E38").Value = "  +5.52%  "
$ws.Range("D39").Value = "'2.630"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "'0.9164"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "'2.048"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").Value = "'105.75"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'0.9962"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "'0.4197"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").Value = "'5.574"
$ws.Range("E45").Value = "  -5.80%  "
$ws.Range("D46").Value = "'7.216"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "'61.87"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").Value = "'0.1227"
$ws.Range("E48").Value = "  -5.37%  "
$ws.Range("D49").Value = "'8.901"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("D51").Value = "'33.54"
$ws.Range("E51").Value = "  -1.02%  "
